$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 197 (pushes existing rows 197-207 down to 198-208)
$ws.Rows.Item(197).Insert()

# Populate the new row 197 with the new weekly data entry
$ws.Cells.Item(197, 1).Value() = 7
$ws.Cells.Item(197, 2).Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(197, 3).Value() = "Ñuble"
$ws.Cells.Item(197, 4).Value() = 45021
$ws.Cells.Item(197, 5).Value() = 16
$ws.Cells.Item(197, 6).Value() = 100112040
$ws.Cells.Item(197, 7).Value() = "Cilantro"
$ws.Cells.Item(197, 8).Value() = "Sin especificar"
$ws.Cells.Item(197, 9).Value() = "Primera"
$ws.Cells.Item(197, 10).Value() = 350
$ws.Cells.Item(197, 11).Value() = 1400
$ws.Cells.Item(197, 12).Value() = 1500
$ws.Cells.Item(197, 13).Value() = 1443
$ws.Cells.Item(197, 14).Value() = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(197, 15).Value() = "Provincia de Diguillín"
$ws.Cells.Item(197, 16).Value() = 1443
$ws.Cells.Item(197, 17).Value() = 1
$ws.Cells.Item(197, 18).Value() = "Hortaliza"
